$d = $word.ActiveDocument

function Toggle-Range($start, $end) {
    $r = $d.Range($start, $end)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

function Get-Text() {
    return $d.Content.Text
}

# =========================================================================
# PHASE 1: text content edits only (insert / replace). Doing text surgery
# first (before any run-splitting) avoids the engine's Range.Text / 
# InsertAfter side effect of re-merging whole paragraphs into a single run.
# =========================================================================

# ---- Edit 1 text: " Obviously the data" -> " Obviously, the data" ----
$full = Get-Text
$idx = $full.IndexOf("Obviously the data says a lot")
if ($idx -lt 0) { throw "anchor not found: Obviously the data says a lot" }
$afterObviously = $idx + ("Obviously").Length
$insPos = $d.Range($afterObviously, $afterObviously)
$insPos.InsertAfter(",")

# ---- Edit 2 text: no textual change (only a run split later) ----

# ---- Edit 3 text: insert "a " before "splendid"; "sub-categories" -> "sub-category" ----
$full = Get-Text
$idx = $full.IndexOf("what if data are packaged and polished using splendid")
if ($idx -lt 0) { throw "anchor not found: what if data are packaged..." }
$usingEnd = $idx + ("what if data are packaged and polished using ").Length
$insPos = $d.Range($usingEnd, $usingEnd)
$insPos.InsertAfter("a ")

$full = Get-Text
$idx = $full.IndexOf("sub-categories")
if ($idx -lt 0) { throw "anchor not found: sub-categories" }
$catEnd = $idx + ("sub-categor").Length
$iesEnd = $catEnd + ("ies").Length
$r = $d.Range($catEnd, $iesEnd)
$r.Text = "y"

# ---- Edit 4 text: insert "a " after "took " ----
$full = Get-Text
$idx = $full.IndexOf("took huge effort to finish a well-designed data visualization")
if ($idx -lt 0) { throw "anchor not found: took huge effort..." }
$tookEnd = $idx + ("took ").Length
$insPos = $d.Range($tookEnd, $tookEnd)
$insPos.InsertAfter("a ")

# ---- Edit 5 text: insert "," after "globe" ----
$full = Get-Text
$idx = $full.IndexOf("globe again and again")
if ($idx -lt 0) { throw "anchor not found: globe again and again" }
$globeEnd = $idx + ("globe").Length
$insPos = $d.Range($globeEnd, $globeEnd)
$insPos.InsertAfter(",")

# ---- Edit 6 text: "users' needs" (curly) -> "user's needs" (straight) ----
$full = Get-Text
$idx = $full.IndexOf([char]0x2019 + " needs and provide ")
if ($idx -lt 0) { throw "anchor not found: users' needs and provide" }
$usersWord = "users"
$searchWindowStart = $idx - 10
if ($searchWindowStart -lt 0) { $searchWindowStart = 0 }
$windowLen = $idx - $searchWindowStart
$window = $full.Substring($searchWindowStart, $windowLen)
$usersPosInWindow = $window.LastIndexOf($usersWord)
if ($usersPosInWindow -lt 0) { throw "could not locate 'users' before apostrophe" }
$usersStart = $searchWindowStart + $usersPosInWindow
$userEnd = $usersStart + ("user").Length
$apostropheOldEnd = $idx + 1
$apRange = $d.Range($userEnd, $apostropheOldEnd)
$apRange.Text = "'s"

Write-Output "Phase1 (text) done"

# =========================================================================
# PHASE 2: run splitting via a harmless Bold-on/off toggle, computed AFTER
# all text mutations above are final (so nothing merges the splits back).
# =========================================================================

# ---- Edit 1 split: " Obviously" | "," | " the data says a lot..." ----
$full = Get-Text
$idx = $full.IndexOf("Obviously, the data says a lot")
if ($idx -lt 0) { throw "split anchor not found: Obviously, the data says a lot" }
$spaceStart = $idx - 1
$obvEnd = $idx + ("Obviously").Length
$commaEnd = $obvEnd + 1
Toggle-Range $spaceStart $obvEnd
Toggle-Range $obvEnd $commaEnd

# ---- Edit 2 split: " that answer" | " to" | " that goal." ----
$full = Get-Text
$idx = $full.IndexOf("that answer to that goal.")
if ($idx -lt 0) { throw "split anchor not found: that answer to that goal." }
$answerStart = $idx - 1
$answerEnd = $answerStart + (" that answer").Length
$toEnd = $answerEnd + (" to").Length
Toggle-Range $answerStart $answerEnd
Toggle-Range $answerEnd $toEnd

# ---- Edit 3 split: "...using " | "a " | "splendid...sub-categor" | "y" ----
$full = Get-Text
$idx = $full.IndexOf("what if data are packaged and polished using a splendid")
if ($idx -lt 0) { throw "split anchor not found: ...using a splendid" }
$usingStart = $idx - 1
$usingEnd = $usingStart + (" what if data are packaged and polished using ").Length
$aEnd = $usingEnd + ("a ").Length
$categorEnd = $aEnd + ("splendid color palette to show the proportion of every sub-categor").Length
$yEnd = $categorEnd + ("y").Length
Toggle-Range $usingStart $usingEnd
Toggle-Range $usingEnd $aEnd
Toggle-Range $aEnd $categorEnd
Toggle-Range $categorEnd $yEnd

# ---- Edit 4 split: "took " | "a " | "huge effort..." ----
$full = Get-Text
$idx = $full.IndexOf("took a huge effort to finish a well-designed data visualization")
if ($idx -lt 0) { throw "split anchor not found: took a huge effort..." }
$tookStart = $idx
$tookEnd = $tookStart + ("took ").Length
$aEnd2 = $tookEnd + ("a ").Length
Toggle-Range $tookStart $tookEnd
Toggle-Range $tookEnd $aEnd2

# ---- Edit 5 split: new "," run between "globe" and " again and again" ----
$full = Get-Text
$idx = $full.IndexOf("globe, again and again")
if ($idx -lt 0) { throw "split anchor not found: globe, again and again" }
$globeEnd = $idx + ("globe").Length
$commaEnd2 = $globeEnd + 1
Toggle-Range $globeEnd $commaEnd2

# ---- Edit 6 split: "...our user" | "'s" | " needs and provide " ----
$full = Get-Text
$idx = $full.IndexOf("user's needs and provide")
if ($idx -lt 0) { throw "split anchor not found: user's needs and provide" }
$userStart2 = $idx
$userEnd2 = $userStart2 + ("user").Length
$apEnd2 = $userEnd2 + ("'s").Length
Toggle-Range $userStart2 $userEnd2
Toggle-Range $userEnd2 $apEnd2

Write-Output "Phase2 (splits) done"

# =========================================================================
# Edit 7: remove w:hint="eastAsia" from the paragraph mark run properties
#   of the final "Ending words" paragraph.
# =========================================================================
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$lastPara.Range.Font.NameAscii = "Times New Roman"

Write-Output "Edit7 done"

Write-Output "All edits complete"
